# Recompute NATMI Lamb1-Itga6 LR-pair stats after raising the number of
# ligand/receptor-expressing cells (E,K) from 1 to 3 per Dr Hou's advice.
# Ligand/receptor/edge average, total, and specificity values updated accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 36.81423266666667
$ws.Range("H2").Value = 110.442698
$ws.Range("I2").Value = 0.13776174071044
$ws.Range("J2").Value = 0.13776174071044
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 192.8285726666667
$ws.Range("N2").Value = 578.485718
$ws.Range("O2").Value = 0.7801188850698786
$ws.Range("P2").Value = 0.7801188850698786
$ws.Range("Q2").Value = 7098.835938931908
$ws.Range("R2").Value = 63889.52345038717
$ws.Range("S2").Value = 0.1074705355683142
$ws.Range("T2").Value = 0.1074705355683142

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 36.81423266666667
$ws.Range("H3").Value = 110.442698
$ws.Range("I3").Value = 0.13776174071044
$ws.Range("J3").Value = 0.13776174071044
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.4209206666666667
$ws.Range("N3").Value = 1.262762
$ws.Range("O3").Value = 0.001702901995496819
$ws.Range("P3").Value = 0.001702901995496819
$ws.Range("Q3").Value = 15.49587135687511
$ws.Range("R3").Value = 139.462842211876
$ws.Range("S3").Value = 0.0002345947431589237
$ws.Range("T3").Value = 0.0002345947431589237

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 36.81423266666667
$ws.Range("H4").Value = 110.442698
$ws.Range("I4").Value = 0.13776174071044
$ws.Range("J4").Value = 0.13776174071044
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 45.70525533333333
$ws.Range("N4").Value = 137.115766
$ws.Range("O4").Value = 0.184907933193646
$ws.Range("P4").Value = 0.184907933193646
$ws.Range("Q4").Value = 1682.603903930741
$ws.Range("R4").Value = 15143.43513537667
$ws.Range("S4").Value = 0.02547323874792643
$ws.Range("T4").Value = 0.02547323874792643

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 36.81423266666667
$ws.Range("H5").Value = 110.442698
$ws.Range("I5").Value = 0.13776174071044
$ws.Range("J5").Value = 0.13776174071044
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.223696
$ws.Range("N5").Value = 24.671088
$ws.Range("O5").Value = 0.0332702797409786
$ws.Range("P5").Value = 0.0332702797409786
$ws.Range("Q5").Value = 302.749057923936
$ws.Range("R5").Value = 2724.741521315425
$ws.Range("S5").Value = 0.0045833716510405
$ws.Range("T5").Value = 0.004583371651040499

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 214.101181
$ws.Range("H6").Value = 642.303543
$ws.Range("I6").Value = 0.8011833806175486
$ws.Range("J6").Value = 0.8011833806175486
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 192.8285726666667
$ws.Range("N6").Value = 578.485718
$ws.Range("O6").Value = 0.7801188850698786
$ws.Range("P6").Value = 0.7801188850698786
$ws.Range("Q6").Value = 41284.82513847765
$ws.Range("R6").Value = 371563.4262462989
$ws.Range("S6").Value = 0.6250182856238782
$ws.Range("T6").Value = 0.6250182856238782

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 214.101181
$ws.Range("H7").Value = 642.303543
$ws.Range("I7").Value = 0.8011833806175486
$ws.Range("J7").Value = 0.8011833806175486
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.4209206666666667
$ws.Range("N7").Value = 1.262762
$ws.Range("O7").Value = 0.001702901995496819
$ws.Range("P7").Value = 0.001702901995496819
$ws.Range("Q7").Value = 90.11961184064067
$ws.Range("R7").Value = 811.0765065657661
$ws.Range("S7").Value = 0.001364336777612511
$ws.Range("T7").Value = 0.001364336777612511

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 214.101181
$ws.Range("H8").Value = 642.303543
$ws.Range("I8").Value = 0.8011833806175486
$ws.Range("J8").Value = 0.8011833806175486
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 45.70525533333333
$ws.Range("N8").Value = 137.115766
$ws.Range("O8").Value = 0.184907933193646
$ws.Range("P8").Value = 0.184907933193646
$ws.Range("Q8").Value = 9785.549144773215
$ws.Range("R8").Value = 88069.94230295894
$ws.Range("S8").Value = 0.1481451630190892
$ws.Range("T8").Value = 0.1481451630190892

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 214.101181
$ws.Range("H9").Value = 642.303543
$ws.Range("I9").Value = 0.8011833806175486
$ws.Range("J9").Value = 0.8011833806175486
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.223696
$ws.Range("N9").Value = 24.671088
$ws.Range("O9").Value = 0.0332702797409786
$ws.Range("P9").Value = 0.0332702797409786
$ws.Range("Q9").Value = 1760.703025784976
$ws.Range("R9").Value = 15846.32723206478
$ws.Range("S9").Value = 0.02665559519696878
$ws.Range("T9").Value = 0.02665559519696877

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.2781493333333334
$ws.Range("H10").Value = 0.8344480000000001
$ws.Range("I10").Value = 0.001040856580779521
$ws.Range("J10").Value = 0.001040856580779521
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 192.8285726666667
$ws.Range("N10").Value = 578.485718
$ws.Range("O10").Value = 0.7801188850698786
$ws.Range("P10").Value = 0.7801188850698786
$ws.Range("Q10").Value = 53.63513893485156
$ws.Range("R10").Value = 482.7162504136641
$ws.Range("S10").Value = 0.0008119918753153661
$ws.Range("T10").Value = 0.0008119918753153661

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.2781493333333334
$ws.Range("H11").Value = 0.8344480000000001
$ws.Range("I11").Value = 0.001040856580779521
$ws.Range("J11").Value = 0.001040856580779521
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.4209206666666667
$ws.Range("N11").Value = 1.262762
$ws.Range("O11").Value = 0.001702901995496819
$ws.Range("P11").Value = 0.001702901995496819
$ws.Range("Q11").Value = 0.1170788028195556
$ws.Range("R11").Value = 1.053709225376
$ws.Range("S11").Value = 0.000001772476748435443
$ws.Range("T11").Value = 0.000001772476748435442

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.2781493333333334
$ws.Range("H12").Value = 0.8344480000000001
$ws.Range("I12").Value = 0.001040856580779521
$ws.Range("J12").Value = 0.001040856580779521
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 45.70525533333333
$ws.Range("N12").Value = 137.115766
$ws.Range("O12").Value = 0.184907933193646
$ws.Range("P12").Value = 0.184907933193646
$ws.Range("Q12").Value = 12.71288630079645
$ws.Range("R12").Value = 114.415976707168
$ws.Range("S12").Value = 0.0001924626391029466
$ws.Range("T12").Value = 0.0001924626391029466

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.2781493333333334
$ws.Range("H13").Value = 0.8344480000000001
$ws.Range("I13").Value = 0.001040856580779521
$ws.Range("J13").Value = 0.001040856580779521
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 8.223696
$ws.Range("N13").Value = 24.671088
$ws.Range("O13").Value = 0.0332702797409786
$ws.Range("P13").Value = 0.0332702797409786
$ws.Range("Q13").Value = 2.287415559936
$ws.Range("R13").Value = 20.586740039424
$ws.Range("S13").Value = 0.00003462958961277317
$ws.Range("T13").Value = 0.00003462958961277316

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 16.037618
$ws.Range("H14").Value = 48.112854
$ws.Range("I14").Value = 0.06001402209123193
$ws.Range("J14").Value = 0.06001402209123194
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 192.8285726666667
$ws.Range("N14").Value = 578.485718
$ws.Range("O14").Value = 0.7801188850698786
$ws.Range("P14").Value = 0.7801188850698786
$ws.Range("Q14").Value = 3092.510987913241
$ws.Range("R14").Value = 27832.59889121917
$ws.Range("S14").Value = 0.04681807200237092
$ws.Range("T14").Value = 0.04681807200237093

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 16.037618
$ws.Range("H15").Value = 48.112854
$ws.Range("I15").Value = 0.06001402209123193
$ws.Range("J15").Value = 0.06001402209123194
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.4209206666666667
$ws.Range("N15").Value = 1.262762
$ws.Range("O15").Value = 0.001702901995496819
$ws.Range("P15").Value = 0.001702901995496819
$ws.Range("Q15").Value = 6.750564860305333
$ws.Range("R15").Value = 60.75508374274801
$ws.Range("S15").Value = 0.000102197997976949
$ws.Range("T15").Value = 0.000102197997976949

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 16.037618
$ws.Range("H16").Value = 48.112854
$ws.Range("I16").Value = 0.06001402209123193
$ws.Range("J16").Value = 0.06001402209123194
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 45.70525533333333
$ws.Range("N16").Value = 137.115766
$ws.Range("O16").Value = 0.184907933193646
$ws.Range("P16").Value = 0.184907933193646
$ws.Range("Q16").Value = 733.0034256284627
$ws.Range("R16").Value = 6597.030830656164
$ws.Range("S16").Value = 0.01109706878752751
$ws.Range("T16").Value = 0.01109706878752751

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 16.037618
$ws.Range("H17").Value = 48.112854
$ws.Range("I17").Value = 0.06001402209123193
$ws.Range("J17").Value = 0.06001402209123194
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 8.223696
$ws.Range("N17").Value = 24.671088
$ws.Range("O17").Value = 0.0332702797409786
$ws.Range("P17").Value = 0.0332702797409786
$ws.Range("Q17").Value = 131.888494996128
$ws.Range("R17").Value = 1186.996454965152
$ws.Range("S17").Value = 0.001996683303356556
$ws.Range("T17").Value = 0.001996683303356556
